# pyFIFOtax test workbook: add a "comment" column (header only, bold like
# the rest of the header row) to every data sheet, and switch the active /
# selected sheet from "currency_conversions" to "sell_orders".

$wb = $excel.ActiveWorkbook

# sheet name -> first empty column after the existing headers
$commentColumns = @{
    "rsu"                  = "G"
    "espp"                 = "G"
    "dividends"            = "F"
    "buy_orders"           = "G"
    "sell_orders"          = "G"
    "currency_conversions" = "F"
}

foreach ($sheetName in $commentColumns.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $col = $commentColumns[$sheetName]
    $headerCell = $ws.Range($col + "1")
    $headerCell.Value = "comment"
    $headerCell.Font.Bold = $true
}

# "sell_orders" becomes the active/selected sheet instead of
# "currency_conversions".
$wb.Worksheets.Item("sell_orders").Activate()
